# Auto-generated edit script: applies numeric cell updates across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (103 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 463.2
$ws.Range("I4").Value = 172
$ws.Range("K4").Value = 172
$ws.Range("M4").Value = -58
$ws.Range("H9").Value = 1425
$ws.Range("J9").Value = 5100
$ws.Range("L9").Value = 5100
$ws.Range("N9").Value = -5438
$ws.Range("H17").Value = 2392
$ws.Range("J17").Value = 2392
$ws.Range("L17").Value = 7176
$ws.Range("N17").Value = -7512
$ws.Range("H41").Value = 278.69232
$ws.Range("I41").Value = 330.375
$ws.Range("J41").Value = 196
$ws.Range("K41").Value = 330.375
$ws.Range("L41").Value = 196
$ws.Range("M41").Value = 109.625
$ws.Range("N41").Value = -1076
$ws.Range("H62").Value = 16781.05
$ws.Range("I62").Value = 16534
$ws.Range("J62").Value = 16914.076
$ws.Range("K62").Value = 16534
$ws.Range("L62").Value = 16914.076
$ws.Range("M62").Value = -15910
$ws.Range("N62").Value = -18162.076
$ws.Range("H65").Value = 16781.05
$ws.Range("I65").Value = 16534
$ws.Range("J65").Value = 16914.076
$ws.Range("K65").Value = 82670
$ws.Range("L65").Value = 84570.38
$ws.Range("M65").Value = -79550
$ws.Range("N65").Value = -90810.38
$ws.Range("H69").Value = 11914.1
$ws.Range("I69").Value = 6315.3335
$ws.Range("K69").Value = 18946.0005
$ws.Range("M69").Value = -18072.0005
$ws.Range("H72").Value = 11914.1
$ws.Range("I72").Value = 6315.3335
$ws.Range("K72").Value = 56838.0015
$ws.Range("M72").Value = -52470.0015
$ws.Range("H74").Value = 5483.5713
$ws.Range("I74").Value = 5398.3335
$ws.Range("J74").Value = 5995
$ws.Range("K74").Value = 5398.3335
$ws.Range("L74").Value = 5995
$ws.Range("M74").Value = -4462.3335
$ws.Range("N74").Value = -7867
$ws.Range("H77").Value = 5483.5713
$ws.Range("I77").Value = 5398.3335
$ws.Range("J77").Value = 5995
$ws.Range("K77").Value = 26991.6675
$ws.Range("L77").Value = 29975
$ws.Range("M77").Value = -22311.6675
$ws.Range("N77").Value = -39335
$ws.Range("H94").Value = 7455.2856
$ws.Range("I94").Value = 7364.8335
$ws.Range("J94").Value = 7998
$ws.Range("K94").Value = 7364.8335
$ws.Range("L94").Value = 7998
$ws.Range("M94").Value = -6913.8335
$ws.Range("N94").Value = -8900
$ws.Range("H97").Value = 2935.2856
$ws.Range("I97").Value = 875
$ws.Range("J97").Value = 3759.4
$ws.Range("K97").Value = 2625
$ws.Range("L97").Value = 11278.2
$ws.Range("M97").Value = -2129
$ws.Range("N97").Value = -12270.2
$ws.Range("H98").Value = 1529.8
$ws.Range("I98").Value = 1238.25
$ws.Range("K98").Value = 1238.25
$ws.Range("M98").Value = 259.75
$ws.Range("H103").Value = 1462.0667
$ws.Range("I103").Value = 1625.25
$ws.Range("J103").Value = 1446.1464
$ws.Range("K103").Value = 4875.75
$ws.Range("L103").Value = 4338.439200000001
$ws.Range("M103").Value = -4289.75
$ws.Range("N103").Value = -5510.439200000001
$ws.Range("H122").Value = 1529.8
$ws.Range("I122").Value = 1238.25
$ws.Range("K122").Value = 3714.75
$ws.Range("M122").Value = -1264.75
$ws.Range("H129").Value = 1213.0714
$ws.Range("J129").Value = 1162.4445
$ws.Range("L129").Value = 3487.3335
$ws.Range("N129").Value = -13487.3335
$ws.Range("H131").Value = 2228.0588
$ws.Range("I131").Value = 1117.3125
$ws.Range("K131").Value = 3351.9375
$ws.Range("M131").Value = 1688.0625
$ws.Range("H135").Value = 1643.125
$ws.Range("I135").Value = 1561
$ws.Range("K135").Value = 14049
$ws.Range("M135").Value = -11514
$ws.Range("H138").Value = 4105.921
$ws.Range("I138").Value = 2779.0952
$ws.Range("J138").Value = 5744.9414
$ws.Range("K138").Value = 8337.285600000001
$ws.Range("L138").Value = 17234.8242
$ws.Range("M138").Value = -3197.285600000001
$ws.Range("N138").Value = -27514.8242

# --- Sheet: ARM (24 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H32").Value = 28374.906
$ws.Range("I32").Value = 28466.365
$ws.Range("K32").Value = 28466.365
$ws.Range("M32").Value = -28179.365
$ws.Range("H74").Value = 41452.58
$ws.Range("I74").Value = 42950.68
$ws.Range("K74").Value = 42950.68
$ws.Range("M74").Value = -42076.68
$ws.Range("H77").Value = 41452.58
$ws.Range("I77").Value = 42950.68
$ws.Range("K77").Value = 214753.4
$ws.Range("M77").Value = -210385.4
$ws.Range("H125").Value = 144999.67
$ws.Range("J125").Value = 144999.67
$ws.Range("L125").Value = 144999.67
$ws.Range("N125").Value = -154839.67
$ws.Range("H132").Value = 32789
$ws.Range("I132").Value = 35661.234
$ws.Range("K132").Value = 106983.702
$ws.Range("M132").Value = -104453.702

# --- Sheet: BSM (14 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2431.5386
$ws.Range("I94").Value = 1793
$ws.Range("J94").Value = 5943.5
$ws.Range("K94").Value = 1793
$ws.Range("L94").Value = 5943.5
$ws.Range("M94").Value = -1342
$ws.Range("N94").Value = -6845.5
$ws.Range("H134").Value = 1716.7255
$ws.Range("I134").Value = 1364.4348
$ws.Range("J134").Value = 4957.8
$ws.Range("K134").Value = 4093.3044
$ws.Range("L134").Value = 14873.4
$ws.Range("M134").Value = -1558.3044
$ws.Range("N134").Value = -19943.4

# --- Sheet: CRP (28 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2211.4814
$ws.Range("I31").Value = 2075.087
$ws.Range("K31").Value = 2075.087
$ws.Range("M31").Value = -1780.087
$ws.Range("H34").Value = 2211.4814
$ws.Range("I34").Value = 2075.087
$ws.Range("K34").Value = 2075.087
$ws.Range("M34").Value = -1873.087
$ws.Range("H58").Value = 52781.85
$ws.Range("J58").Value = 3279.2222
$ws.Range("L58").Value = 3279.2222
$ws.Range("N58").Value = -3685.2222
$ws.Range("H99").Value = 14852.25
$ws.Range("I99").Value = 21679.2
$ws.Range("K99").Value = 21679.2
$ws.Range("M99").Value = -20181.2
$ws.Range("H122").Value = 2022.9565
$ws.Range("I122").Value = 2076.476
$ws.Range("K122").Value = 6229.428
$ws.Range("M122").Value = -3779.428
$ws.Range("H126").Value = 14852.25
$ws.Range("I126").Value = 21679.2
$ws.Range("K126").Value = 65037.60000000001
$ws.Range("M126").Value = -62567.60000000001
$ws.Range("H136").Value = 52781.85
$ws.Range("J136").Value = 3279.2222
$ws.Range("L136").Value = 9837.6666
$ws.Range("N136").Value = -14937.6666

# --- Sheet: CUL (31 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 613.3333
$ws.Range("I5").Value = 613.3333
$ws.Range("K5").Value = 1839.9999
$ws.Range("M5").Value = -1727.9999
$ws.Range("H33").Value = 46.285713
$ws.Range("I33").Value = 24
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 144
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = 139
$ws.Range("N33").Value = -866
$ws.Range("H134").Value = 5145.75
$ws.Range("I134").Value = 5113.2
$ws.Range("K134").Value = 15339.6
$ws.Range("M134").Value = -10269.6
$ws.Range("H135").Value = 613.3333
$ws.Range("I135").Value = 613.3333
$ws.Range("K135").Value = 5519.9997
$ws.Range("M135").Value = -2984.9997
$ws.Range("H137").Value = 4820.0835
$ws.Range("I137").Value = 2820.5
$ws.Range("K137").Value = 8461.5
$ws.Range("M137").Value = -3361.5
$ws.Range("H138").Value = 455.42856
$ws.Range("I138").Value = 455.42856
$ws.Range("K138").Value = 1366.28568
$ws.Range("M138").Value = 3773.71432
$ws.Range("H140").Value = 2056.875
$ws.Range("I140").Value = 1910.7333
$ws.Range("K140").Value = 5732.199900000001
$ws.Range("M140").Value = -552.1999000000005

# --- Sheet: GSM (18 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 18982.334
$ws.Range("I102").Value = 1854.3334
$ws.Range("J102").Value = 36110.332
$ws.Range("K102").Value = 1854.3334
$ws.Range("L102").Value = 36110.332
$ws.Range("M102").Value = -232.3334
$ws.Range("N102").Value = -39354.332
$ws.Range("H122").Value = 2980.5
$ws.Range("I122").Value = 2721.125
$ws.Range("K122").Value = 8163.375
$ws.Range("M122").Value = -5713.375
$ws.Range("H132").Value = 27973.77
$ws.Range("I132").Value = 36549.207
$ws.Range("J132").Value = 3105
$ws.Range("K132").Value = 109647.621
$ws.Range("L132").Value = 9315
$ws.Range("M132").Value = -107117.621
$ws.Range("N132").Value = -14375

# --- Sheet: LTW (13 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3174.25
$ws.Range("H27").Value = 3174.25
$ws.Range("H40").Value = 9418.467000000001
$ws.Range("I40").Value = 10501.692
$ws.Range("J40").Value = 2377.5
$ws.Range("K40").Value = 10501.692
$ws.Range("L40").Value = 2377.5
$ws.Range("M40").Value = -10365.692
$ws.Range("N40").Value = -2649.5
$ws.Range("H132").Value = 55228.543
$ws.Range("J132").Value = 6089.1
$ws.Range("L132").Value = 18267.3
$ws.Range("N132").Value = -23327.3

# --- Sheet: WVR (15 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 24038
$ws.Range("I49").Value = 24038
$ws.Range("K49").Value = 24038
$ws.Range("M49").Value = -23808
$ws.Range("H122").Value = 47984.184
$ws.Range("I122").Value = 2757
$ws.Range("J122").Value = 251506.5
$ws.Range("K122").Value = 8271
$ws.Range("L122").Value = 754519.5
$ws.Range("M122").Value = -5821
$ws.Range("N122").Value = -759419.5
$ws.Range("H132").Value = 30560.205
$ws.Range("I132").Value = 31941.838
$ws.Range("K132").Value = 95825.514
$ws.Range("M132").Value = -93295.514

